$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Delete()

$null = $ws.Hyperlinks.Add("E2", "https://site-B.com")
